# Add data for 2021-11-17
# - Rename the "through" date from 2021-11-08 to 2021-11-09 (sheet name +
#   header label in column B, row 1).
# - Update/insert carjacking counts for the new day across several
#   neighborhood rows / month columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet name + running-total header label
$ws.Name = "Through 2021-11-09"
$ws.Range("B1").Value = "November 2021 (through November 09)"

# Row 2 - North Lawndale
$ws.Range("M2").Value = 6
$ws.Range("AI2").Value = 1
$ws.Range("BE2").Value = 1
$ws.Range("BP2").Value = 1

# Row 3 - Garfield Park
$ws.Range("X3").Value = 2

# Row 4 - Austin
$ws.Range("B4").Value = 3

# Row 6 - West Town
$ws.Range("B6").Value = 6

# Row 7 - Englewood
$ws.Range("B7").Value = 3

# Row 8 - South Shore
$ws.Range("B8").Value = 1
$ws.Range("M8").Value = 2

# Row 9 - Auburn Gresham
$ws.Range("BP9").Value = 2

# Row 13 - Loop
$ws.Range("B13").Value = 1

# Row 17 - West Loop
$ws.Range("AT17").Value = 2

# Row 20 - Lake View
$ws.Range("M20").Value = 2

# Row 22 - Little Village
$ws.Range("B22").Value = 2

# Row 25 - Chinatown
$ws.Range("B25").Value = 1

# Row 27 - River North
$ws.Range("M27").Value = 1

# Row 30 - Edgewater
$ws.Range("AI30").Value = 1
$ws.Range("BB30").Value = 1
$ws.Range("BE30").Value = 1

# Row 32 - Little Italy, UIC
$ws.Range("M32").Value = 3

# Row 35 - Hyde Park
$ws.Range("AT35").Value = 1

# Row 39 - West Elsdon
$ws.Range("M39").Value = 1

# Row 43 - Ashburn
$ws.Range("B43").Value = 3

# Row 47 - Roseland
$ws.Range("AI47").Value = 1
$ws.Range("AT47").Value = 2

# Row 64 - Bridgeport
$ws.Range("B64").Value = 2

# Row 72 - Gage Park
$ws.Range("M72").Value = 1
$ws.Range("BE72").Value = 1
